$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, pushing existing data down.
$ws.Rows.Item(1).Insert()

$headers = @("Workbook", "# sheets", "Sheetname", "Rows", "Columns", "Workbook Size", "Useful(1-5)", "Description")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
